$d = $word.ActiveDocument

# Each (old, new) pair is unique in the document, so a simple
# Find/Replace over the whole content is safe and order-independent.
$replacements = @(
    @('2025-10-28 Tuesday', '2025-10-29 Wednesday'),
    @('542÷2=271, 0', '103÷7=14, 5'),
    @('101÷2=50, 1', '401÷9=44, 5'),
    @('627÷5=125, 2', '454÷6=75, 4'),
    @('321÷9=35, 6', '219÷4=54, 3'),
    @('407÷6=67, 5', '539÷8=67, 3'),
    @('691÷5=138, 1', '588÷9=65, 3'),
    @('912÷7=130, 2', '991÷4=247, 3'),
    @('787÷7=112, 3', '892÷7=127, 3'),
    @('830÷6=138, 2', '702÷5=140, 2'),
    @('341÷6=56, 5', '712÷8=89, 0'),
    @('558÷8=69, 6', '260÷3=86, 2'),
    @('445÷3=148, 1', '377÷2=188, 1'),
    @('158÷9=17, 5', '299÷5=59, 4'),
    @('862÷2=431, 0', '867÷4=216, 3'),
    @('501÷8=62, 5', '790÷2=395, 0'),
    @('608÷3=202, 2', '296÷2=148, 0'),
    @('147÷5=29, 2', '658÷8=82, 2'),
    @('115÷6=19, 1', '963÷6=160, 3'),
    @('652÷8=81, 4', '770÷7=110, 0'),
    @('725÷4=181, 1', '869÷4=217, 1'),
    @('682÷2=341, 0', '684÷4=171, 0'),
    @('617÷7=88, 1', '692÷4=173, 0'),
    @('826÷8=103, 2', '664÷6=110, 4'),
    @('799÷6=133, 1', '713÷4=178, 1'),
    @('267÷2=133, 1', '934÷6=155, 4')
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

Write-Output 'done'
